$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain text before assigning values that look numeric
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.879.28"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "1.826.51"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "310.85"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").Value = "0.4571"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D8").Value = "0.3682"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "0.07167"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("D10").Value = "0.8727"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "0.07753"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "19.58"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "1.816.89"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "5.314"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "6.377"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").Value = "86.90"
$ws.Range("E16").Value = "  -5.44%  "
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "0.000008707"
$ws.Range("E18").Value = "  -3.95%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "26.917.48"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "14.46"
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").Value = "5.000"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").Value = "2.053.01"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").Value = "10.44"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").Value = "2.005"
$ws.Range("E25").Value = "  +4.53%  "
$ws.Range("D26").Value = "151.55"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").Value = "18.17"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").Value = "1.955"
$ws.Range("E28").Value = "  -5.60%  "
$ws.Range("D29").Value = "113.56"
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").Value = "4.902"
$ws.Range("E30").Value = "  -4.22%  "
$ws.Range("D31").Value = "0.08792"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "3.044"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").Value = "0.7487"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("D34").Value = "4.480"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  -3.60%  "
$ws.Range("D36").Value = "2.542"
$ws.Range("E36").Value = "  -3.57%  "
$ws.Range("D37").Value = "1.081"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "0.01943"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").Value = "0.05125"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").Value = "2.908"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").Value = "6.922"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("D42").Value = "0.4966"
$ws.Range("E42").Value = "  -3.59%  "
$ws.Range("D43").Value = "0.1593"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "8.296"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").Value = "0.4689"
$ws.Range("E45").Value = "  -3.18%  "
$ws.Range("D46").Value = "1.007"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "10.14"
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("D48").Value = "101.50"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").Value = "1.610"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("D50").Value = "0.06101"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "64.48"
$ws.Range("E51").Value = "  -2.03%  "

# Restore original (default) style so the cells do not carry an explicit text number format
$ws.Range("D2:D51").Style = "Normal"
